# Generate Report for Handoff
# Refresh the "Latest Handoff Date/Datetime" column for every row whose
# handoff was re-run as part of this handoff generation pass (rows 7 and
# 10-16 on each sheet -- the rows currently carrying the most recent batch
# of handoff timestamps), so they all show the single, newer timestamp for
# that locale.

$wb = $excel.ActiveWorkbook

# --- Overview sheet: column D = "Latest Handoff Date" ---
$ws = $wb.Worksheets.Item("Overview")
$newDate = "2016-22-20 00:22:52"
foreach ($r in 7,10,11,12,13,14,15,16) {
    $ws.Cells.Item($r, 4).Value = $newDate
}

# --- zh-cn sheet: column E = "Latest Handoff Datetime" ---
$ws = $wb.Worksheets.Item("zh-cn")
$newDateTimeZhCn = "2016-03-20 00:22:49"
foreach ($r in 7,10,11,12,13,14,15,16) {
    $ws.Cells.Item($r, 5).Value = $newDateTimeZhCn
}

# --- de-de sheet: column E = "Latest Handoff Datetime" ---
$ws = $wb.Worksheets.Item("de-de")
$newDateTimeDeDe = "2016-03-20 00:22:52"
foreach ($r in 7,10,11,12,13,14,15,16) {
    $ws.Cells.Item($r, 5).Value = $newDateTimeDeDe
}
